$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The document already contains an identical "<del>...</del>" markup
# pattern elsewhere (e.g. "Aulcuns <del>po</del> ont invente..."). Use
# it as a formatting reference so the freshly-created runs end up with
# exactly the same rPr (Courier New, red a91111, size 9pt/18 half-pts)
# that the rest of the document's <del> markers use.
# ------------------------------------------------------------------
$refRange = $d.Content
$refRange.Find.Execute("Aulcuns <del>po</del>") | Out-Null
$refStart = $refRange.Start

$lenAulcunsSp = "Aulcuns ".Length
$lenDelOpen   = "<del>".Length
$lenPo        = "po".Length
$lenDelClose  = "</del>".Length

$refDelOpenStart  = $refStart + $lenAulcunsSp
$refDelOpenEnd    = $refDelOpenStart + $lenDelOpen
$refPoEnd         = $refDelOpenEnd + $lenPo
$refDelCloseEnd   = $refPoEnd + $lenDelClose

$refDelOpen  = $d.Range($refDelOpenStart, $refDelOpenEnd)   # "<del>"
$refPo       = $d.Range($refDelOpenEnd,   $refPoEnd)        # "po"
$refDelClose = $d.Range($refPoEnd,        $refDelCloseEnd)  # "</del>"

# ------------------------------------------------------------------
# Find the target text to edit: "Aulcuns ne mectent" -> split it into
# "Aulcuns " + "<del>" + "po" + "</del>" + " ne mectent"
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Aulcuns ne mectent") | Out-Null
$start = $target.Start

# Insert the raw characters first (plain text insert keeps the single
# run's current formatting; we re-format the relevant sub-ranges next).
$insertAt = $start + $lenAulcunsSp
$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.Text = "<del>po</del> "

# Re-apply formatting to each newly created segment by copying the
# FormattedText from the matching reference run, so rFonts/color/size
# match exactly.
$delOpenStart = $insertAt
$delOpenEnd   = $delOpenStart + $lenDelOpen
$poEnd        = $delOpenEnd + $lenPo
$delCloseEnd  = $poEnd + $lenDelClose

$delOpenRange = $d.Range($delOpenStart, $delOpenEnd)   # "<del>"
$delOpenRange.FormattedText = $refDelOpen.FormattedText

$poRange = $d.Range($delOpenEnd, $poEnd)               # "po"
$poRange.FormattedText = $refPo.FormattedText

$delCloseRange = $d.Range($poEnd, $delCloseEnd)        # "</del>"
$delCloseRange.FormattedText = $refDelClose.FormattedText
